$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.185.75"
$ws.Range("E2").Value = "  -4.04%  "
$ws.Range("D3").Value = "3.296.50"
$ws.Range("E3").Value = "  -4.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.05"
$ws.Range("E5").Value = "  -3.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.68"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.297.87"
$ws.Range("E8").Value = "  -4.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.86"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "3.869.34"
$ws.Range("E13").Value = "  -4.19%  "
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.17"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").Value = "3.306.16"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "60.241.64"
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.22"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.52"
$ws.Range("E21").Value = "  -5.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.57"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.81"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.546"
$ws.Range("E24").Value = "  -3.70%  "
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "3.454.92"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("E27").Value = "  -9.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.174"
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  -6.97%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.60"
$ws.Range("E33").Value = "  -5.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.50"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.27"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.17"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.75"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.71"
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.51"
$ws.Range("E39").Value = "  -7.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.98"
$ws.Range("E40").Value = "  -14.23%  "
$ws.Range("D41").Value = "3.338.08"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0733"
$ws.Range("E42").Value = "  -5.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.86"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.749"
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.18"
$ws.Range("E45").Value = "  -4.52%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  -5.82%  "
$ws.Range("D48").Value = "2.356.20"
$ws.Range("E48").Value = "  -8.11%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.53"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.32"
$ws.Range("E51").Value = "  -6.04%  "
